$wb = $excel.ActiveWorkbook

# --- TestReports sheet: remove the "TestResultExcelFilePath" column (H) ---
$wsTestReports = $wb.Worksheets.Item("TestReports")
$wsTestReports.Columns("H").Select()
$wsTestReports.Columns("H").Delete()

# --- ProcessPayrollForNIWeekly sheet: remove the "TestResultExcelFilePath" column (H) ---
$wsProcess = $wb.Worksheets.Item("ProcessPayrollForNIWeekly")
$wsProcess.Activate()
$wsProcess.Columns("H").Delete()
$wsProcess.Range("M7").Select()
